# Update user operation interface:
#  - Add a new "requesterId" field row to the first table (nodeDraw), row 11
#    (row 11 was previously blank, just like the row above the table
#    separators elsewhere in the sheet, so no rows need to be shifted)
#  - Add a brand-new "用户表userData" table describing the user model
#    (userName / userId / userPassword / userLevel) in rows 44-48, right
#    after the last existing table
#  - Refresh view state (zoom + selection) to match the edited workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add "requesterId" to the first table (节点绘制表nodeDraw) ---
$ws.Range("A11").Value = "requesterId"

# --- 2. Add a new table block describing 用户表userData ---
# Header row (style matches the other table headers: bold black for A:C,
# bold red for the table-name cell in column E)
$ws.Range("A44").Value = "字段"
$ws.Range("B44").Value = "类型"
$ws.Range("C44").Value = "注释"
$ws.Range("E44").Value = "用户表userData（用户包括申请人，1级处理人，2级处理人等）"
$ws.Range("A44:C44").Font.Bold = $true
$ws.Range("E44").Font.Bold = $true
$ws.Range("E44").Font.Color = 255

# Data rows
$ws.Range("A45").Value = "userName"
$ws.Range("B45").Value = "varchar(100)"
$ws.Range("C45").Value = "用户名"

$ws.Range("A46").Value = "userId"
$ws.Range("B46").Value = "varchar(64)"
$ws.Range("C46").Value = "用户id"

$ws.Range("A47").Value = "userPassword"
$ws.Range("B47").Value = "varchar(100)"
$ws.Range("C47").Value = "用户密码"

$ws.Range("A48").Value = "userLevel"
$ws.Range("B48").Value = "varchar(100)"
$ws.Range("C48").Value = "用户权限（权限从0开始逐渐增加）"

# --- 3. View state: zoom to 130% and move selection to E13 ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("E13").Select() | Out-Null
